$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 11321.6
$ws.Range("I9").Value = 2029.5
$ws.Range("K9").Value = 2029.5
$ws.Range("M9").Value = -1860.5

$ws.Range("H41").Value = 1183.4445
$ws.Range("I41").Value = 2069.5
$ws.Range("J41").Value = 474.6
$ws.Range("K41").Value = 2069.5
$ws.Range("L41").Value = 474.6
$ws.Range("M41").Value = -1629.5
$ws.Range("N41").Value = -1354.6

$ws.Range("H86").Value = 1553.5294
$ws.Range("I86").Value = 1481.0769
$ws.Range("K86").Value = 1481.0769
$ws.Range("M86").Value = -358.0769

$ws.Range("H89").Value = 1553.5294
$ws.Range("I89").Value = 1481.0769
$ws.Range("K89").Value = 7405.3845
$ws.Range("M89").Value = -1789.3845

$ws.Range("H115").Value = 1055
$ws.Range("I115").Value = 1055
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3165
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1598
$ws.Range("N115").ClearContents()

$ws.Range("H137").Value = 100456.414
$ws.Range("I137").Value = 2103.0613
$ws.Range("J137").Value = 354104.53
$ws.Range("K137").Value = 6309.1839
$ws.Range("L137").Value = 1062313.59
$ws.Range("M137").Value = -3759.1839
$ws.Range("N137").Value = -1067413.59

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1507.1562
$ws.Range("I2").Value = 1753.174
$ws.Range("J2").Value = 878.44446
$ws.Range("K2").Value = 1753.174
$ws.Range("L2").Value = 878.44446
$ws.Range("M2").Value = -1640.174
$ws.Range("N2").Value = -1104.44446

$ws.Range("H24").Value = 42009.668
$ws.Range("J24").Value = 42009.668
$ws.Range("L24").Value = 42009.668
$ws.Range("N24").Value = -42757.668

$ws.Range("H63").Value = 1891.3462
$ws.Range("I63").Value = 1870.909
$ws.Range("K63").Value = 1870.909
$ws.Range("M63").Value = -1184.909

$ws.Range("H66").Value = 1891.3462
$ws.Range("I66").Value = 1870.909
$ws.Range("K66").Value = 9354.545
$ws.Range("M66").Value = -5922.545

$ws.Range("H93").Value = 75407
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H97").Value = 832.30554
$ws.Range("I97").Value = 535.069
$ws.Range("K97").Value = 535.069
$ws.Range("M97").Value = -39.06899999999996

$ws.Range("H100").Value = 42009.668
$ws.Range("J100").Value = 42009.668
$ws.Range("L100").Value = 42009.668
$ws.Range("N100").Value = -44173.668

$ws.Range("H116").Value = 1507.1562
$ws.Range("I116").Value = 1753.174
$ws.Range("J116").Value = 878.44446
$ws.Range("K116").Value = 1753.174
$ws.Range("L116").Value = 878.44446
$ws.Range("M116").Value = 540.826
$ws.Range("N116").Value = -5466.44446

$ws.Range("H122").Value = 30337.879
$ws.Range("I122").Value = 2837.92
$ws.Range("J122").Value = 116275.25
$ws.Range("K122").Value = 8513.76
$ws.Range("L122").Value = 348825.75
$ws.Range("M122").Value = -6063.76
$ws.Range("N122").Value = -353725.75

$ws.Range("H132").Value = 2376.7346
$ws.Range("I132").Value = 2284.8809
$ws.Range("K132").Value = 6854.6427
$ws.Range("M132").Value = -4324.6427

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1507.1562
$ws.Range("I3").Value = 1753.174
$ws.Range("J3").Value = 878.44446
$ws.Range("K3").Value = 1753.174
$ws.Range("L3").Value = 878.44446
$ws.Range("M3").Value = -1639.174
$ws.Range("N3").Value = -1106.44446

$ws.Range("H22").Value = 206.375
$ws.Range("I22").Value = 213
$ws.Range("K22").Value = 213
$ws.Range("M22").Value = -40

$ws.Range("H86").Value = 17577732
$ws.Range("I86").Value = 24414376
$ws.Range("K86").Value = 24414376
$ws.Range("M86").Value = -24413253

$ws.Range("H89").Value = 17577732
$ws.Range("I89").Value = 24414376
$ws.Range("K89").Value = 122071880
$ws.Range("M89").Value = -122066264

$ws.Range("H94").Value = 32653.285
$ws.Range("I94").Value = 493.72726
$ws.Range("J94").Value = 150571.67
$ws.Range("K94").Value = 493.72726
$ws.Range("L94").Value = 150571.67
$ws.Range("M94").Value = -42.72726
$ws.Range("N94").Value = -151473.67

$ws.Range("H99").Value = 2011.7142
$ws.Range("J99").Value = 1915.6
$ws.Range("L99").Value = 1915.6
$ws.Range("N99").Value = -4911.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 651.1818
$ws.Range("I22").Value = 261.6
$ws.Range("J22").Value = 975.8333
$ws.Range("K22").Value = 261.6
$ws.Range("L22").Value = 975.8333
$ws.Range("M22").Value = 88.39999999999998
$ws.Range("N22").Value = -1675.8333

$ws.Range("H31").Value = 2333.368
$ws.Range("I31").Value = 1903.738
$ws.Range("J31").Value = 2734.3555
$ws.Range("K31").Value = 1903.738
$ws.Range("L31").Value = 2734.3555
$ws.Range("M31").Value = -1608.738
$ws.Range("N31").Value = -3324.3555

$ws.Range("H34").Value = 2333.368
$ws.Range("I34").Value = 1903.738
$ws.Range("J34").Value = 2734.3555
$ws.Range("K34").Value = 1903.738
$ws.Range("L34").Value = 2734.3555
$ws.Range("M34").Value = -1701.738
$ws.Range("N34").Value = -3138.3555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1251.4814
$ws.Range("J12").Value = 1578.55
$ws.Range("L12").Value = 4735.65
$ws.Range("N12").Value = -5081.65

$ws.Range("H75").Value = 682.1429000000001
$ws.Range("I75").Value = 736
$ws.Range("J75").Value = 610.3333
$ws.Range("K75").Value = 2208
$ws.Range("L75").Value = 1830.9999
$ws.Range("M75").Value = -1210
$ws.Range("N75").Value = -3826.9999

$ws.Range("H78").Value = 682.1429000000001
$ws.Range("I78").Value = 736
$ws.Range("J78").Value = 610.3333
$ws.Range("K78").Value = 6624
$ws.Range("L78").Value = 5492.9997
$ws.Range("M78").Value = -1632
$ws.Range("N78").Value = -15476.9997

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H138").Value = 5881.3
$ws.Range("I138").Value = 3156
$ws.Range("K138").Value = 9468
$ws.Range("M138").Value = -4328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 290.5
$ws.Range("I2").Value = 130.2
$ws.Range("J2").Value = 450.8
$ws.Range("K2").Value = 130.2
$ws.Range("L2").Value = 450.8
$ws.Range("M2").Value = -17.19999999999999
$ws.Range("N2").Value = -676.8

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H122").Value = 130254.36
$ws.Range("I122").Value = 164978.27
$ws.Range("J122").Value = 2933.3333
$ws.Range("K122").Value = 494934.8099999999
$ws.Range("L122").Value = 8799.999899999999
$ws.Range("M122").Value = -492484.8099999999
$ws.Range("N122").Value = -13699.9999

$ws.Range("H126").Value = 56545.06
$ws.Range("I126").Value = 85482.37
$ws.Range("K126").Value = 256447.11
$ws.Range("M126").Value = -253977.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 41530.188
$ws.Range("I2").Value = 60555.11
$ws.Range("K2").Value = 60555.11
$ws.Range("M2").Value = -60443.11

$ws.Range("H22").Value = 2784.93
$ws.Range("I22").Value = 4121.8696
$ws.Range("K22").Value = 4121.8696
$ws.Range("M22").Value = -3826.8696

$ws.Range("H27").Value = 2784.93
$ws.Range("I27").Value = 4121.8696
$ws.Range("K27").Value = 4121.8696
$ws.Range("M27").Value = -4014.8696

$ws.Range("H39").Value = 20298
$ws.Range("I39").Value = 19980.666
$ws.Range("J39").Value = 21250
$ws.Range("K39").Value = 19980.666
$ws.Range("L39").Value = 21250
$ws.Range("M39").Value = -19520.666
$ws.Range("N39").Value = -22170

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H58").Value = 13033.333
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4740

$ws.Range("H60").Value = 40666
$ws.Range("J60").Value = 40666
$ws.Range("L60").Value = 40666
$ws.Range("N60").Value = -41684

$ws.Range("H82").Value = 1944.3043
$ws.Range("I82").Value = 1985.9474
$ws.Range("K82").Value = 1985.9474
$ws.Range("M82").Value = -1624.9474

$ws.Range("H85").Value = 1944.3043
$ws.Range("I85").Value = 1985.9474
$ws.Range("K85").Value = 1985.9474
$ws.Range("M85").Value = -737.9474

$ws.Range("H93").Value = 1160.9286
$ws.Range("I93").Value = 1037.4546
$ws.Range("K93").Value = 1037.4546
$ws.Range("M93").Value = 210.5454

$ws.Range("H94").Value = 29982.166
$ws.Range("J94").Value = 29982.166
$ws.Range("L94").Value = 29982.166
$ws.Range("N94").Value = -31334.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 9512500
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H122").Value = 3684.6667
$ws.Range("I122").Value = 3380.5881
$ws.Range("J122").Value = 4977
$ws.Range("K122").Value = 10141.7643
$ws.Range("L122").Value = 14931
$ws.Range("M122").Value = -7691.764299999999
$ws.Range("N122").Value = -19831

$ws.Range("H126").Value = 2973.7856
$ws.Range("I126").Value = 2973.7856
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8921.356800000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6451.356800000001
$ws.Range("N126").ClearContents()
